# Commit: "base de datos con docker"
#
# The canonical diff shows the final slide of the deck (sldId 261 /
# r:id rId7 -> ppt/slides/slide6.xml, the "QUE ES DOCKER?" slide) being
# removed entirely: it disappears from <p:sldIdLst> in presentation.xml,
# and the slide6.xml part (plus its relationship entries) is dropped
# from the package.
#
# Deleting the slide through the PowerPoint object model takes care of
# all of that in one shot: it removes the <p:sldId> entry, drops the
# slide part, and cleans up the now-unused relationship / content-type
# entries that pointed at it.

$p = $ppt.ActivePresentation

# The slide to remove is the last one in the deck (6 slides -> 5).
$lastIndex = $p.Slides.Count
$slide = $p.Slides.Item($lastIndex)
$slide.Delete()
